# Case_0_250 / res_line / pl_mw.xlsx: recomputed line active-power-loss
# results for the "380 kV" case (commit: "case with 380 kV done").
# Columns B1:O1 are the line indices (header row, unchanged); rows 2-25 are
# per-scenario results. Only columns C,D,E,F,G,H,J,K,L,O change per row -
# B, I, M, N stay 0 and column A (scenario index) is untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 3).Value = 0.4268873391042973
$ws.Cells.Item(2, 4).Value = 0.2046583537790383
$ws.Cells.Item(2, 5).Value = 0.1864598340150678
$ws.Cells.Item(2, 6).Value = 1.512273373307259
$ws.Cells.Item(2, 7).Value = 0.8654205197321545
$ws.Cells.Item(2, 8).Value = 0.9421265453390291
$ws.Cells.Item(2, 10).Value = 0.2196388092749473
$ws.Cells.Item(2, 11).Value = 2.062959816201669
$ws.Cells.Item(2, 12).Value = 0.1694129256959869
$ws.Cells.Item(2, 15).Value = 3.644308653382978

# Row 3
$ws.Cells.Item(3, 3).Value = 0.4216639451440898
$ws.Cells.Item(3, 4).Value = 0.1997997793834969
$ws.Cells.Item(3, 5).Value = 0.1850631748933154
$ws.Cells.Item(3, 6).Value = 1.526682210855391
$ws.Cells.Item(3, 7).Value = 0.8777951587670429
$ws.Cells.Item(3, 8).Value = 0.9533504526584053
$ws.Cells.Item(3, 10).Value = 0.2199948724294671
$ws.Cells.Item(3, 11).Value = 1.847978277075299
$ws.Cells.Item(3, 12).Value = 0.1687248952590608
$ws.Cells.Item(3, 15).Value = 3.693707738598761

# Row 4
$ws.Cells.Item(4, 3).Value = 0.4186434002496213
$ws.Cells.Item(4, 4).Value = 0.1968656607029828
$ws.Cells.Item(4, 5).Value = 0.1842738000879507
$ws.Cells.Item(4, 6).Value = 1.536488585419974
$ws.Cells.Item(4, 7).Value = 0.8860941790756485
$ws.Cells.Item(4, 8).Value = 0.9607471702700963
$ws.Cells.Item(4, 10).Value = 0.2203112764400785
$ws.Cells.Item(4, 11).Value = 1.715574019572614
$ws.Cells.Item(4, 12).Value = 0.1683519026339724
$ws.Cells.Item(4, 15).Value = 3.72657131162066

# Row 5
$ws.Cells.Item(5, 3).Value = 0.417459563553777
$ws.Cells.Item(5, 4).Value = 0.195682441333048
$ws.Cells.Item(5, 5).Value = 0.1839693216577594
$ws.Cells.Item(5, 6).Value = 1.540725823342704
$ws.Cells.Item(5, 7).Value = 0.8896520333237277
$ws.Cells.Item(5, 8).Value = 0.9638883843885324
$ws.Cells.Item(5, 10).Value = 0.2204648263980715
$ws.Cells.Item(5, 11).Value = 1.661520352166235
$ws.Cells.Item(5, 12).Value = 0.1682123869000449
$ws.Cells.Item(5, 15).Value = 3.740599481397481

# Row 6
$ws.Cells.Item(6, 3).Value = 0.4172658344437536
$ws.Cells.Item(6, 4).Value = 0.1954867246206788
$ws.Cells.Item(6, 5).Value = 0.1839198035673846
$ws.Cells.Item(6, 6).Value = 1.541443968226467
$ws.Cells.Item(6, 7).Value = 0.8902534289758748
$ws.Cells.Item(6, 8).Value = 0.9644176493139724
$ws.Cells.Item(6, 10).Value = 0.2204918105467435
$ws.Cells.Item(6, 11).Value = 1.652538987528146
$ws.Cells.Item(6, 12).Value = 0.1681899755996703
$ws.Cells.Item(6, 15).Value = 3.742967238069113

# Row 7
$ws.Cells.Item(7, 3).Value = 0.4186272439000049
$ws.Cells.Item(7, 4).Value = 0.1968496527913857
$ws.Cells.Item(7, 5).Value = 0.1842696240785138
$ws.Cells.Item(7, 6).Value = 1.536544754462611
$ws.Cells.Item(7, 7).Value = 0.8861414496120759
$ws.Cells.Item(7, 8).Value = 0.9607890196569215
$ws.Cells.Item(7, 10).Value = 0.2203132475788863
$ws.Cells.Item(7, 11).Value = 1.714845423913289
$ws.Cells.Item(7, 12).Value = 0.1683499704783884
$ws.Cells.Item(7, 15).Value = 3.72675792633872

# Row 8
$ws.Cells.Item(8, 3).Value = 0.4250476495788007
$ws.Cells.Item(8, 4).Value = 0.2029730036558988
$ws.Cells.Item(8, 5).Value = 0.1859641467609912
$ws.Cells.Item(8, 6).Value = 1.51704234316454
$ws.Cells.Item(8, 7).Value = 0.8695416224171737
$ws.Cells.Item(8, 8).Value = 0.9458916705175895
$ws.Cells.Item(8, 10).Value = 0.2197413003493764
$ws.Cells.Item(8, 11).Value = 1.988920451549234
$ws.Cells.Item(8, 12).Value = 0.1691654533271247
$ws.Cells.Item(8, 15).Value = 3.660815406883486

# Row 9
$ws.Cells.Item(9, 3).Value = 0.4391144431215821
$ws.Cells.Item(9, 4).Value = 0.2153656587447301
$ws.Cells.Item(9, 5).Value = 0.1898260457953462
$ws.Cells.Item(9, 6).Value = 1.486417339562891
$ws.Cells.Item(9, 7).Value = 0.8425654274084522
$ws.Cells.Item(9, 8).Value = 0.920688261704683
$ws.Cells.Item(9, 10).Value = 0.2193946389922772
$ws.Cells.Item(9, 11).Value = 2.523024635944694
$ws.Cells.Item(9, 12).Value = 0.1711553156196857
$ws.Cells.Item(9, 15).Value = 3.551629731263574

# Row 10
$ws.Cells.Item(10, 3).Value = 0.4503449362916569
$ws.Cells.Item(10, 4).Value = 0.2247001078808211
$ws.Cells.Item(10, 5).Value = 0.1929897126018858
$ws.Cells.Item(10, 6).Value = 1.468573121818423
$ws.Cells.Item(10, 7).Value = 0.8261649723789759
$ws.Cells.Item(10, 8).Value = 0.9046182557261488
$ws.Cells.Item(10, 10).Value = 0.2196114236165485
$ws.Cells.Item(10, 11).Value = 2.9132302695769
$ws.Cells.Item(10, 12).Value = 0.1728534068756034
$ws.Cells.Item(10, 15).Value = 3.483727949963935

# Row 11
$ws.Cells.Item(11, 3).Value = 0.4556476257817224
$ws.Cells.Item(11, 4).Value = 0.2289954547205184
$ws.Cells.Item(11, 5).Value = 0.194499369835949
$ws.Cells.Item(11, 6).Value = 1.461468764094413
$ws.Cells.Item(11, 7).Value = 0.8194505874845248
$ws.Cells.Item(11, 8).Value = 0.8978393972849119
$ws.Cells.Item(11, 10).Value = 0.2198122447814086
$ws.Cells.Item(11, 11).Value = 3.090237336870814
$ws.Cells.Item(11, 12).Value = 0.1736767776778336
$ws.Cells.Item(11, 15).Value = 3.45552219929769

# Row 12
$ws.Cells.Item(12, 3).Value = 0.4576833938143068
$ws.Cells.Item(12, 4).Value = 0.2306289414007665
$ws.Cells.Item(12, 5).Value = 0.1950811278878639
$ws.Cells.Item(12, 6).Value = 1.458924403298994
$ws.Cells.Item(12, 7).Value = 0.8170156879842096
$ws.Cells.Item(12, 8).Value = 0.8953488978507806
$ws.Cells.Item(12, 10).Value = 0.2199029680323363
$ws.Cells.Item(12, 11).Value = 3.15719021058095
$ws.Cells.Item(12, 12).Value = 0.1739958464720743
$ws.Cells.Item(12, 15).Value = 3.445228107744271

# Row 13
$ws.Cells.Item(13, 3).Value = 0.4572437223412464
$ws.Cells.Item(13, 4).Value = 0.2302768344549406
$ws.Cells.Item(13, 5).Value = 0.1949553882465374
$ws.Cells.Item(13, 6).Value = 1.459465884133166
$ws.Cells.Item(13, 7).Value = 0.81753529146944
$ws.Cells.Item(13, 8).Value = 0.8958818674716724
$ws.Cells.Item(13, 10).Value = 0.2198827767314526
$ws.Cells.Item(13, 11).Value = 3.142774149114246
$ws.Cells.Item(13, 12).Value = 0.1739268063772599
$ws.Cells.Item(13, 15).Value = 3.447427905214369

# Row 14
$ws.Cells.Item(14, 3).Value = 0.4558145541148519
$ws.Cells.Item(14, 4).Value = 0.2291297044230163
$ws.Cells.Item(14, 5).Value = 0.1945470296309217
$ws.Cells.Item(14, 6).Value = 1.461256512542128
$ws.Cells.Item(14, 7).Value = 0.8192481063483825
$ws.Cells.Item(14, 8).Value = 0.8976329686066151
$ws.Cells.Item(14, 10).Value = 0.2198194146116919
$ws.Cells.Item(14, 11).Value = 3.095747133450516
$ws.Cells.Item(14, 12).Value = 0.1737028820822104
$ws.Cells.Item(14, 15).Value = 3.454667540222687

# Row 15
$ws.Cells.Item(15, 3).Value = 0.4549427580667782
$ws.Cells.Item(15, 4).Value = 0.2284279532893123
$ws.Cells.Item(15, 5).Value = 0.1942982098144022
$ws.Cells.Item(15, 6).Value = 1.462372332054329
$ws.Cells.Item(15, 7).Value = 0.8203112909055221
$ws.Cells.Item(15, 8).Value = 0.8987155345570699
$ws.Cells.Item(15, 10).Value = 0.219782514214991
$ws.Cells.Item(15, 11).Value = 3.066931729916178
$ws.Cells.Item(15, 12).Value = 0.1735666682732244
$ws.Cells.Item(15, 15).Value = 3.45915242953032

# Row 16
$ws.Cells.Item(16, 3).Value = 0.450002283733312
$ws.Cells.Item(16, 4).Value = 0.2244203735598518
$ws.Cells.Item(16, 5).Value = 0.1928924668695373
$ws.Cells.Item(16, 6).Value = 1.469057813271974
$ws.Cells.Item(16, 7).Value = 0.826618822523244
$ws.Cells.Item(16, 8).Value = 0.9050719716282458
$ws.Cells.Item(16, 10).Value = 0.2196003543651628
$ws.Cells.Item(16, 11).Value = 2.901652044884486
$ws.Cells.Item(16, 12).Value = 0.1728006181999504
$ws.Cells.Item(16, 15).Value = 3.485625337504302

# Row 17
$ws.Cells.Item(17, 3).Value = 0.4470210311858125
$ws.Cells.Item(17, 4).Value = 0.2219743309258888
$ws.Cells.Item(17, 5).Value = 0.1920481055641403
$ws.Cells.Item(17, 6).Value = 1.473418772393082
$ws.Cells.Item(17, 7).Value = 0.8306796962290477
$ws.Cells.Item(17, 8).Value = 0.9091076203909836
$ws.Cells.Item(17, 10).Value = 0.2195147678191205
$ws.Cells.Item(17, 11).Value = 2.800127503123065
$ws.Cells.Item(17, 12).Value = 0.172343677262873
$ws.Cells.Item(17, 15).Value = 3.502553520591903

# Row 18
$ws.Cells.Item(18, 3).Value = 0.4453245501138667
$ws.Cells.Item(18, 4).Value = 0.2205720568510543
$ws.Cells.Item(18, 5).Value = 0.1915690891964204
$ws.Cells.Item(18, 6).Value = 1.476022421811486
$ws.Cells.Item(18, 7).Value = 0.8330856226573999
$ws.Cells.Item(18, 8).Value = 0.9114788381386987
$ws.Cells.Item(18, 10).Value = 0.2194751619025013
$ws.Cells.Item(18, 11).Value = 2.741686501675929
$ws.Cells.Item(18, 12).Value = 0.1720856504980688
$ws.Cells.Item(18, 15).Value = 3.512542652022404

# Row 19
$ws.Cells.Item(19, 3).Value = 0.4447532901773741
$ws.Cells.Item(19, 4).Value = 0.2200980690321899
$ws.Cells.Item(19, 5).Value = 0.1914080441630368
$ws.Cells.Item(19, 6).Value = 1.476920340738907
$ws.Cells.Item(19, 7).Value = 0.8339122766224989
$ws.Cells.Item(19, 8).Value = 0.912290280074977
$ws.Cells.Item(19, 10).Value = 0.2194634050903659
$ws.Cells.Item(19, 11).Value = 2.721891473873825
$ws.Cells.Item(19, 12).Value = 0.1719991116936157
$ws.Cells.Item(19, 15).Value = 3.515968135305229

# Row 20
$ws.Cells.Item(20, 3).Value = 0.4473365020364781
$ws.Cells.Item(20, 4).Value = 0.2222342385664007
$ws.Cells.Item(20, 5).Value = 0.1921373026490087
$ws.Cells.Item(20, 6).Value = 1.472944671751478
$ws.Cells.Item(20, 7).Value = 0.8302401390964178
$ws.Cells.Item(20, 8).Value = 0.9086728412783955
$ws.Cells.Item(20, 10).Value = 0.2195228830424867
$ws.Cells.Item(20, 11).Value = 2.810939833891666
$ws.Cells.Item(20, 12).Value = 0.1723918236032063
$ws.Cells.Item(20, 15).Value = 3.500725347855507

# Row 21
$ws.Cells.Item(21, 3).Value = 0.4562335832101212
$ws.Cells.Item(21, 4).Value = 0.2294664571595746
$ws.Cells.Item(21, 5).Value = 0.1946667012068986
$ws.Cells.Item(21, 6).Value = 1.460726599937836
$ws.Cells.Item(21, 7).Value = 0.8187420853357992
$ws.Cells.Item(21, 8).Value = 0.8971165507150261
$ws.Cells.Item(21, 10).Value = 0.2198376274439227
$ws.Cells.Item(21, 11).Value = 3.10956219874447
$ws.Cells.Item(21, 12).Value = 0.1737684569467177
$ws.Cells.Item(21, 15).Value = 3.452530578062607

# Row 22
$ws.Cells.Item(22, 3).Value = 0.4622100573367334
$ws.Cells.Item(22, 4).Value = 0.2342334791032528
$ws.Cells.Item(22, 5).Value = 0.1963785601910928
$ws.Cells.Item(22, 6).Value = 1.453591936139837
$ws.Cells.Item(22, 7).Value = 0.8118553314061359
$ws.Cells.Item(22, 8).Value = 0.890009827359421
$ws.Cells.Item(22, 10).Value = 0.2201288684260732
$ws.Cells.Item(22, 11).Value = 3.304285156741912
$ws.Cells.Item(22, 12).Value = 0.1747105574163044
$ws.Cells.Item(22, 15).Value = 3.423287613946314

# Row 23
$ws.Cells.Item(23, 3).Value = 0.4590055447752661
$ws.Cells.Item(23, 4).Value = 0.231685578125763
$ws.Cells.Item(23, 5).Value = 0.1954595504629495
$ws.Cells.Item(23, 6).Value = 1.45732194000766
$ws.Cells.Item(23, 7).Value = 0.8154733443924016
$ws.Cells.Item(23, 8).Value = 0.8937619847709612
$ws.Cells.Item(23, 10).Value = 0.2199656076354799
$ws.Cells.Item(23, 11).Value = 3.200399801180367
$ws.Cells.Item(23, 12).Value = 0.1742038761285443
$ws.Cells.Item(23, 15).Value = 3.438688467052401

# Row 24
$ws.Cells.Item(24, 3).Value = 0.4471938232056232
$ws.Cells.Item(24, 4).Value = 0.2221167218907283
$ws.Cells.Item(24, 5).Value = 0.1920969566515396
$ws.Cells.Item(24, 6).Value = 1.473158712149747
$ws.Cells.Item(24, 7).Value = 0.8304386408821429
$ws.Cells.Item(24, 8).Value = 0.9088692458401084
$ws.Cells.Item(24, 10).Value = 0.2195191842487887
$ws.Cells.Item(24, 11).Value = 2.806051806826531
$ws.Cells.Item(24, 12).Value = 0.1723700420788958
$ws.Cells.Item(24, 15).Value = 3.501551064606801

# Row 25
$ws.Cells.Item(25, 3).Value = 0.4351514044125793
$ws.Cells.Item(25, 4).Value = 0.2119724445998799
$ws.Cells.Item(25, 5).Value = 0.1887238331296146
$ws.Cells.Item(25, 6).Value = 1.4938853764234
$ws.Cells.Item(25, 7).Value = 0.8492642142711375
$ws.Cells.Item(25, 8).Value = 0.9270769016972693
$ws.Cells.Item(25, 10).Value = 0.2194055564675637
$ws.Cells.Item(25, 11).Value = 2.378910776246585
$ws.Cells.Item(25, 12).Value = 0.1705753899276772
$ws.Cells.Item(25, 15).Value = 3.579007732502546
